$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Simplify-Rhs($rhs) {
    $termPattern = '(\d+ \* )?\(\*([A-Z])\)\[i - (\d+)\]'
    $termMatches = [regex]::Matches($rhs, $termPattern)
    $parts = New-Object System.Collections.ArrayList
    foreach ($tm in $termMatches) {
        $coefStr = $tm.Groups[1].Value
        $var = $tm.Groups[2].Value
        $idx = $tm.Groups[3].Value
        if ($coefStr -eq '') {
            $coef = 1
        } else {
            $trimmed = $coefStr.Trim()
            $parts2 = $trimmed.Split(' ')
            $coef = [int]$parts2[0]
        }
        if ($coef -ne 0) {
            [void]$parts.Add("(*$var)[i - $idx]")
        }
    }

    if ($parts.Count -eq 1) {
        $newRhs = $parts[0]
    } else {
        $expr = $parts[0]
        for ($i = 1; $i -lt $parts.Count; $i++) {
            $expr = "($expr + $($parts[$i]))"
        }
        if ($expr.StartsWith('(') -and $expr.EndsWith(')')) {
            $expr = $expr.Substring(1, $expr.Length - 2)
        }
        $newRhs = $expr
    }
    return $newRhs
}

function Simplify-Code($code) {
    $stmtPattern = '(\(\*[A-Z]\)\[i\] = )(.*?)(;)'
    $stmtMatches = [regex]::Matches($code, $stmtPattern)
    $result = ""
    $lastEnd = 0
    foreach ($m in $stmtMatches) {
        $start = $m.Index
        $len = $m.Length
        $before = $code.Substring($lastEnd, $start - $lastEnd)
        $result = $result + $before
        $prefix = $m.Groups[1].Value
        $rhs = $m.Groups[2].Value
        $suffix = $m.Groups[3].Value
        $newRhs = Simplify-Rhs $rhs
        $result = $result + $prefix + $newRhs + $suffix
        $lastEnd = $start + $len
    }
    $result = $result + $code.Substring($lastEnd)
    return $result
}

for ($row = 2; $row -le 136; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    $code = $cell.Value2
    $newCode = Simplify-Code $code
    $cell.Value2 = $newCode
}

Write-Host "Done"
